$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 currently holds the blank "just below the table" spacer formatting
# (complete with its distinctive top border). Since the new diary entry will
# push that spacer role down to row 26, move that formatting there first.
$ws.Range("A22:G22").Copy()
$ws.Range("A26:G26").PasteSpecial(-4122)

# Copy the formatting of the last filled diary row (row 21) down onto the
# four rows that will receive the new diary entry (rows 22-25).
$ws.Range("A21:G21").Copy()
$ws.Range("A22:G25").PasteSpecial(-4122)

# The new entries use centered alignment for the Goal/Achievements/Reflection/Mood columns.
$ws.Range("D22:G25").HorizontalAlignment = -4108

# Row 22: 1/30/2020
$ws.Range("A22").Value = "1/30/2020"
$ws.Range("B22").Value = "17:00 - 19:50"
$ws.Range("C22").Value = "None"
$ws.Range("D22").Value = "Looking forward to discuss about our findings as a team and possibly learn the significance of the UML diagrams in detail and the idea behind printing it. Also, wonder who the speaker is gonna be!"
$ws.Range("E22").Value = "As expected, we discussed the details and thought process behind choosing the features that we explored. One of the teams bothered to print out the entire UML diagram. We were very confident while answering the questions and also expressed that we did not use the UML diagram much for the project. We learnt about the essence and the professor gave us tips on effective group work. I was very happy to see the Call graphs and sequence diagrams! Finally, we had a talk with Consuelo Lopez who was in various roles."
$ws.Range("F22").Value = "We realized that a lot of teams found the use of UML rather ambiguous. After learning the details of the arrows, I gained some key understanding of the concepts. Since, the professor advised on more face to face meetings and voicing out, we are planning to implement that during out team meeting for next assignment.We also learnt about a situation where the efficiency was reduced by 1000x, was feeling human. It was good to know that experts also make mistakes and we should always explore the choices and not be reserved. Also, taught us the importance of modeling. After viewing the call graphs, I felt that this would be a better use for the assignment. Listening to Consuelo Lopez, reassured me that I will always have the way to explore variety of fields in Software and our purpose is greater than merely coding. "
$ws.Range("G22").Value = "Productive! Learnt a lot about visualization. Great to hear the talk."

# Row 23: 2/4/2020
$ws.Range("A23").Value = "2/4/2020"
$ws.Range("B23").Value = "13:30 - 15:00"
$ws.Range("C23").Value = "Group"
$ws.Range("D23").Value = "To find two features to build the assignment on!"
$ws.Range("E23").Value = "We figured out the essential features. It was rather easy to identify which are the essential features. "
$ws.Range("F23").Value = "Eventhough it was easy to identify the core features as the app has an interactive UI. It was rather difficult to find how we are going to make the package. We settled on the features and identified the key classes using the technique in previous assignment."
$ws.Range("G23").Value = "Feels neutral!"

# Row 24: 2/5/2020
$ws.Range("A24").Value = "2/5/2020"
$ws.Range("B24").Value = "21:00 - 23:00"
$ws.Range("C24").Value = "Group"
$ws.Range("D24").Value = "To add everything we found to a report and also use the call graphs related to explain further."
$ws.Range("E24").Value = "We finished up the write up using several points we found during our discussion and used call graphs to identify some end function and explain the feature. "
$ws.Range("F24").Value = "We are satisfied with the write up. I think the Call graphs were more useful than UML diagrams because we couldn’t see many relations or the UML was dense. We abandoned UML altogether. Call graphs were easier to follow and make observations to describe in report. "
$ws.Range("G24").Value = "Tired. "

# Row 25: 2/6/2020
$ws.Range("A25").Value = "2/6/2020"
$ws.Range("B25").Value = "8:00 - 10:00"
$ws.Range("C25").Value = "Group"
$ws.Range("D25").Value = "Run over the report and add any missing code or correct grammar, make sure the packet is understandable and print diagrams. "
$ws.Range("E25").Value = "We are satisfied with the report and have done our best to explain with related supporting materials. We printed these call graphs for class. "
$ws.Range("F25").Value = "Relatively easy work as it was just proofreading for today and understanding call graphs. "
$ws.Range("G25").Value = "Satisfied!"

# Row heights for the new entries and the following blank spacer row.
$ws.Rows.Item(22).RowHeight = 316.8
$ws.Rows.Item(23).RowHeight = 107.15
$ws.Rows.Item(24).RowHeight = 101
$ws.Rows.Item(25).RowHeight = 68.5
$ws.Rows.Item(26).RowHeight = 33.85
